$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (A) was left blank as a header for the category/class
# column. Fill it in with the label "classe" to match the other header
# cells in row 1 (years 2020-2039).
$ws.Range("A1").Value = "classe"
